$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "snapshot" column N, mirroring column M (the most recent
#     price-check column) as the starting point for this run's values ---
$ws.Range("M1:M204").Copy($ws.Range("N1"))

# Stamp the new column's header with this run's timestamp.
$ws.Range("N1").Value = "2026-01-28 05:20:04"

# --- Scraper text-cleanup fixes on column B (product name) ---
# Row 9: stray trailing price text left over from a bad scrape; strip it.
$ws.Range("B9").Value = "Samsung Galaxy A17 5G Gris (4 Go / 128 Go)"

# Row 26: this one was previously cleaned but should carry the trailing
# price text like the other still-dirty rows.
$ws.Range("B26").Value = "Apple iPhone 15 128 Go Rose659€00"

# Row 45: stray trailing price text; strip it.
$ws.Range("B45").Value = "Apple iPhone 17 256 Go Blanc"

# Row 91: add back the trailing price text to match the raw scrape.
$ws.Range("B91").Value = "Apple iPhone 17 Pro Max 512 Go Orange Cosmique1 729€00"
